# This edit permutes the per-row data (columns D and K:T) among rows 2-19.
# Columns A,B,C,E,F,G,H,I,J are identical across all rows and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: source row (current/before) -> destination row (after)
$mapping = @{
    2  = 19
    3  = 4
    4  = 13
    5  = 2
    6  = 3
    7  = 10
    8  = 11
    9  = 17
    10 = 18
    11 = 14
    12 = 8
    13 = 9
    14 = 15
    15 = 16
    16 = 5
    17 = 6
    18 = 7
    19 = 12
}

# Columns that carry the per-row data which gets permuted.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot current values for every row/column involved, before any writes.
$snapshot = @{}
foreach ($srcRow in $mapping.Keys) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$srcRow").Value2
    }
    $snapshot[$srcRow] = $rowData
}

# Now write the snapshotted values into their destination rows.
foreach ($srcRow in $mapping.Keys) {
    $destRow = $mapping[$srcRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}
